$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.000158
$ws.Range("E2").Value = 0.300338
$ws.Range("G2").Value = 0.300497
$ws.Range("E3").Value = 0.296978
$ws.Range("G3").Value = 0.297013
$ws.Range("D4").Value = 0.000033
$ws.Range("E4").Value = 0.298187
$ws.Range("G4").Value = 0.29822
$ws.Range("D5").Value = 0.000035
$ws.Range("E5").Value = 0.300219
$ws.Range("G5").Value = 0.300253
$ws.Range("D6").Value = 0.000036
$ws.Range("E6").Value = 0.310343
$ws.Range("G6").Value = 0.310379
$ws.Range("E7").Value = 0.305961
$ws.Range("G7").Value = 0.305997
$ws.Range("D8").Value = 0.000037
$ws.Range("E8").Value = 0.30436
$ws.Range("G8").Value = 0.304397
$ws.Range("E9").Value = 0.301549
$ws.Range("G9").Value = 0.301584
$ws.Range("D10").Value = 0.000034
$ws.Range("E10").Value = 0.295757
$ws.Range("G10").Value = 0.295791
$ws.Range("E11").Value = 0.295803
$ws.Range("G11").Value = 0.295838
$ws.Range("D12").Value = 0.000048
$ws.Range("E12").Value = 0.162024
$ws.Range("G12").Value = 0.162072
$ws.Range("D13").Value = 0.00005
$ws.Range("E13").Value = 0.160357
$ws.Range("G13").Value = 0.160407
$ws.Range("D14").Value = 0.000054
$ws.Range("E14").Value = 0.159372
$ws.Range("G14").Value = 0.159426
$ws.Range("D15").Value = 0.000051
$ws.Range("E15").Value = 0.162482
$ws.Range("G15").Value = 0.162532
$ws.Range("E16").Value = 0.163849
$ws.Range("G16").Value = 0.163906
$ws.Range("D17").Value = 0.000051
$ws.Range("E17").Value = 0.161625
$ws.Range("G17").Value = 0.161677
$ws.Range("D18").Value = 0.000055
$ws.Range("E18").Value = 0.16016
$ws.Range("G18").Value = 0.160215
$ws.Range("D19").Value = 0.000071
$ws.Range("E19").Value = 0.16072
$ws.Range("G19").Value = 0.160791
$ws.Range("D20").Value = 0.00006999999999999999
$ws.Range("E20").Value = 0.160076
$ws.Range("G20").Value = 0.160146
$ws.Range("D21").Value = 0.000054
$ws.Range("E21").Value = 0.161942
$ws.Range("G21").Value = 0.161996
$ws.Range("D22").Value = 0.00009000000000000001
$ws.Range("E22").Value = 0.10532
$ws.Range("G22").Value = 0.10541
$ws.Range("D23").Value = 0.000091
$ws.Range("E23").Value = 0.107003
$ws.Range("G23").Value = 0.107094
$ws.Range("D24").Value = 0.00008500000000000001
$ws.Range("E24").Value = 0.107134
$ws.Range("G24").Value = 0.10722
$ws.Range("D25").Value = 0.000082
$ws.Range("E25").Value = 0.106185
$ws.Range("G25").Value = 0.106267
$ws.Range("D26").Value = 0.000087
$ws.Range("E26").Value = 0.108144
$ws.Range("G26").Value = 0.108231
$ws.Range("E27").Value = 0.112646
$ws.Range("G27").Value = 0.112732
$ws.Range("D28").Value = 0.000097
$ws.Range("E28").Value = 0.108461
$ws.Range("G28").Value = 0.108558
$ws.Range("D29").Value = 0.000084
$ws.Range("E29").Value = 0.104025
$ws.Range("G29").Value = 0.104109
$ws.Range("D30").Value = 0.000088
$ws.Range("E30").Value = 0.109175
$ws.Range("G30").Value = 0.109263
$ws.Range("D31").Value = 0.00008500000000000001
$ws.Range("E31").Value = 0.113916
$ws.Range("G31").Value = 0.114001
